# Generate Report for Handoff
# Updates the localization-status report to reflect that the zh-cn / de-de
# translations are now "Ready for handoff" (previously "In Translation"),
# refreshes the handoff timestamps, and widens the Status columns to fit
# the longer status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -----------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Refresh handoff timestamps --------------------------------------------
# zh-cn handoff datetime
$wsZhCn.Range("H2").Value = "2016-08-30 09:26:08"
# de-de handoff datetime (mirrored on the Overview sheet as the "de-de" column)
$wsDeDe.Range("H2").Value = "2016-08-30 09:26:19"
$wsOverview.Range("G2").Value = "2016-08-30 09:26:19"

# --- Widen the Status columns so the new, longer text fits -----------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
